# Added more test cases
$wb = $excel.ActiveWorkbook

$wsPostorder      = $wb.Worksheets.Item(1)  # postorderTraverse
$wsNumberOfNodes   = $wb.Worksheets.Item(2)  # getNumberOfNodes
$wsHeight         = $wb.Worksheets.Item(3)  # getHeight

# --- getNumberOfNodes: add a 3rd test tree (Tree 1, Node B/C/D) in columns I:K ---
$wsNumberOfNodes.Range("I1").Value = "Tree 1, Node B"
$wsNumberOfNodes.Range("J1").Value = "Tree 1, Node C"
$wsNumberOfNodes.Range("K1").Value = "Tree 1, Node D"

$wsNumberOfNodes.Range("I2").Value = 1
$wsNumberOfNodes.Range("J2").Value = 6
$wsNumberOfNodes.Range("K2").Value = 2

$wsNumberOfNodes.Range("I3").Value = 1
$wsNumberOfNodes.Range("J3").Value = 6
$wsNumberOfNodes.Range("K3").Value = 2

$wsNumberOfNodes.Columns.Item(10).ColumnWidth = 18.5
$wsNumberOfNodes.Columns.Item(11).ColumnWidth = 13.16666666

# --- getHeight: add a 3rd test tree (Tree 1, Node B/C/D) in columns F:H ---
$wsHeight.Range("F1").Value = "Tree 1, Node B"
$wsHeight.Range("G1").Value = "Tree 1, Node C"
$wsHeight.Range("H1").Value = "Tree 1, Node D"

$wsHeight.Range("F2").Value = 1
$wsHeight.Range("G2").Value = 3
$wsHeight.Range("H2").Value = 2

$wsHeight.Range("F3").Value = 1
$wsHeight.Range("G3").Value = 3
$wsHeight.Range("H3").Value = 2

$wsHeight.Columns.Item(8).ColumnWidth = 16.5

# --- Selections on each sheet, matching the saved workbook state ---
$wsPostorder.Range("E11").Select()
$wsNumberOfNodes.Range("I1:K3").Select()
$wsHeight.Range("H6").Select()

# getHeight is the sheet that is active/selected when the workbook is saved
$wsHeight.Activate()
